$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "<m><pa>Eschervis" -> "<pa>Eschervis"
#   (drop the leading "<m>" from the opening-tag run right before "Eschervis")
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("<m><pa>Eschervis", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$openTag = $d.Range($r1.Start, $r1.Start + 7)
$openTag.Text = "<pa>"

# ---------------------------------------------------------------------------
# Edit 2: remove the whole "</m>" run that used to close the <head> math tag
#   ("...Eschervis</pa> racine</m></head>" -> "...Eschervis</pa> racine</head>")
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("</m></head>", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$closeTag = $d.Range($r2.Start, $r2.Start + 4)
$closeTag.Delete()

# ---------------------------------------------------------------------------
# Edit 3: wrap "en lieu fort humide o..ou fontayne on le puysse souvent
#         arroser" with <env> / </env> tags, styled like the other inline
#         tags (blue Courier New, sz 18/18).
#
# A fresh donor Range (whose FormattedText already carries the exact tag
# formatting: rFonts ascii/eastAsia/hAnsi/cs=Courier New, color 0000ff,
# sz/szCs 18) is re-fetched right before each use -- a FormattedText object
# captured before an intervening document mutation can go stale.
# ---------------------------------------------------------------------------

# 3a. Insert "<env>" right after "Ils veulent estre plantés "
$donorA = $d.Content
$donorA.Find.Execute("<oc>", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$envTagFormatA = $donorA.FormattedText

$r3 = $d.Content
$r3.Find.Execute("Ils veulent estre plantés ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$splitPoint = $r3.End
$insPoint = $d.Range($splitPoint, $splitPoint)
$insPoint.FormattedText = $envTagFormatA
$newRun = $d.Range($splitPoint, $splitPoint + 4)
$newRun.Text = "<env>"

# 3b. Insert "</env>" right after "...souvent arroser" (before ", car")
$donorB = $d.Content
$donorB.Find.Execute("<oc>", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$envTagFormatB = $donorB.FormattedText

$r4 = $d.Content
$r4.Find.Execute("ou fontayne on le puysse souvent arroser", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$splitPoint2 = $r4.End
$insPoint2 = $d.Range($splitPoint2, $splitPoint2)
$insPoint2.FormattedText = $envTagFormatB
$newRun2 = $d.Range($splitPoint2, $splitPoint2 + 4)
$newRun2.Text = "</env>"
